# Generate Report for Handoff
# Fill in the "Latest Handoff Datetime" for rows that previously showed a
# placeholder (duplicated) datetime value, with their own distinct
# handoff datetime, for both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhRows = @(4, 6, 7, 8, 9, 10)
$deRows = @(4, 6, 7, 8, 9, 10)

$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $zhRows) {
    $wsZh.Range("D$r").Value = "2016-02-17 05:07:55"
}

$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $deRows) {
    $wsDe.Range("D$r").Value = "2016-02-17 05:08:06"
}
